$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H40").Value = 1910.7
$ws.Range("I40").Value = 1580
$ws.Range("K40").Value = 1580
$ws.Range("M40").Value = -1405

$ws.Range("H116").Value = 17303284
$ws.Range("I116").Value = 34593324
$ws.Range("J116").Value = 13244.25
$ws.Range("K116").Value = 34593324
$ws.Range("L116").Value = 13244.25
$ws.Range("M116").Value = -34589882
$ws.Range("N116").Value = -20128.25

$ws.Range("H129").Value = 1156.9546
$ws.Range("J129").Value = 1247.579
$ws.Range("L129").Value = 3742.737
$ws.Range("N129").Value = -13742.737

$ws.Range("H137").Value = 1284.5518
$ws.Range("I137").Value = 900
$ws.Range("J137").Value = 1298.2858
$ws.Range("K137").Value = 2700
$ws.Range("L137").Value = 3894.8574
$ws.Range("M137").Value = -150
$ws.Range("N137").Value = -8994.857400000001

$ws.Range("H138").Value = 7248931.5
$ws.Range("J138").Value = 10419112
$ws.Range("L138").Value = 31257336
$ws.Range("N138").Value = -31267616

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 2826.2886
$ws.Range("I32").Value = 2980.5293
$ws.Range("K32").Value = 2980.5293
$ws.Range("M32").Value = -2693.5293

$ws.Range("H135").Value = 29800
$ws.Range("J135").Value = 29800
$ws.Range("L135").Value = 29800
$ws.Range("N135").Value = -39940

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H20").Value = 1677.0588
$ws.Range("I20").Value = 1251
$ws.Range("J20").Value = 2285.7144
$ws.Range("K20").Value = 1251
$ws.Range("L20").Value = 2285.7144
$ws.Range("M20").Value = -1004
$ws.Range("N20").Value = -2779.7144

$ws.Range("H82").Value = 19974.223
$ws.Range("I82").Value = 4947.5
$ws.Range("J82").Value = 50027.668
$ws.Range("K82").Value = 4947.5
$ws.Range("L82").Value = 50027.668
$ws.Range("M82").Value = -4564.5
$ws.Range("N82").Value = -50793.668

$ws.Range("H85").Value = 19974.223
$ws.Range("I85").Value = 4947.5
$ws.Range("J85").Value = 50027.668
$ws.Range("K85").Value = 4947.5
$ws.Range("L85").Value = 50027.668
$ws.Range("M85").Value = -3621.5
$ws.Range("N85").Value = -52679.668

$ws.Range("H134").Value = 2056.3845
$ws.Range("I134").Value = 1415.129
$ws.Range("J134").Value = 4541.25
$ws.Range("K134").Value = 4245.387
$ws.Range("L134").Value = 13623.75
$ws.Range("M134").Value = -1710.387
$ws.Range("N134").Value = -18693.75

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H64").Value = 25000
$ws.Range("J64").Value = 25000
$ws.Range("L64").Value = 25000
$ws.Range("N64").Value = -25496

$ws.Range("H67").Value = 25000
$ws.Range("J67").Value = 25000
$ws.Range("L67").Value = 25000
$ws.Range("N67").Value = -26716

$ws.Range("H110").Value = 0
$ws.Range("J110").Value = 0
$ws.Range("L110").Value = 0
$ws.Range("N110").ClearContents()

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H2").Value = 43478416
$ws.Range("I2").Value = 179.42105
$ws.Range("J2").Value = 250000050
$ws.Range("K2").Value = 1076.5263
$ws.Range("L2").Value = 1500000300
$ws.Range("M2").Value = -963.5263
$ws.Range("N2").Value = -1500000526

$ws.Range("H107").Value = 304913.4
$ws.Range("J107").Value = 324535.56
$ws.Range("L107").Value = 973606.6799999999
$ws.Range("N107").Value = -977446.6799999999

$ws.Range("H122").Value = 738.2963
$ws.Range("I122").Value = 561.64703
$ws.Range("J122").Value = 1038.6
$ws.Range("K122").Value = 5054.82327
$ws.Range("L122").Value = 9347.4
$ws.Range("M122").Value = -2604.82327
$ws.Range("N122").Value = -14247.4

$ws.Range("H131").Value = 935.95
$ws.Range("I131").Value = 515
$ws.Range("J131").Value = 962.81915
$ws.Range("K131").Value = 1545
$ws.Range("L131").Value = 2888.45745
$ws.Range("M131").Value = 3495
$ws.Range("N131").Value = -12968.45745

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H70").Value = 5555.8823
$ws.Range("I70").Value = 5587.879
$ws.Range("J70").Value = 4500
$ws.Range("K70").Value = 5587.879
$ws.Range("L70").Value = 4500
$ws.Range("M70").Value = -5317.879
$ws.Range("N70").Value = -5040

$ws.Range("H73").Value = 5555.8823
$ws.Range("I73").Value = 5587.879
$ws.Range("J73").Value = 4500
$ws.Range("K73").Value = 5587.879
$ws.Range("L73").Value = 4500
$ws.Range("M73").Value = -4651.879
$ws.Range("N73").Value = -6372

$ws.Range("H119").Value = 0
$ws.Range("J119").Value = 0
$ws.Range("L119").Value = 0
$ws.Range("N119").ClearContents()

$ws.Range("H122").Value = 1390863.9
$ws.Range("I122").Value = 5556305.5
$ws.Range("J122").Value = 2383.3333
$ws.Range("K122").Value = 16668916.5
$ws.Range("L122").Value = 7149.999899999999
$ws.Range("M122").Value = -16666466.5
$ws.Range("N122").Value = -12049.9999

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H92").Value = 0
$ws.Range("J92").Value = 0
$ws.Range("L92").Value = 0
$ws.Range("N92").ClearContents()

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H32").Value = 0
$ws.Range("J32").Value = 0
$ws.Range("L32").Value = 0
$ws.Range("N32").ClearContents()

$ws.Range("H46").Value = 134842.9
$ws.Range("J46").Value = 134842.9
$ws.Range("L46").Value = 134842.9
$ws.Range("N46").Value = -135304.9

$ws.Range("H70").Value = 0
$ws.Range("J70").Value = 0
$ws.Range("L70").Value = 0
$ws.Range("N70").ClearContents()

$ws.Range("H73").Value = 0
$ws.Range("J73").Value = 0
$ws.Range("L73").Value = 0
$ws.Range("N73").ClearContents()

$ws.Range("H113").Value = 4600.5
$ws.Range("I113").Value = 5634
$ws.Range("J113").Value = 1500
$ws.Range("K113").Value = 16902
$ws.Range("L113").Value = 4500
$ws.Range("M113").Value = -14732
$ws.Range("N113").Value = -8840

$ws.Range("H122").Value = 1000004
$ws.Range("I122").Value = 1000004
$ws.Range("J122").Value = 0
$ws.Range("K122").Value = 3000012
$ws.Range("L122").Value = 0
$ws.Range("M122").Value = -2997562
$ws.Range("N122").ClearContents()

$ws.Range("H126").Value = 101711.6
$ws.Range("I126").Value = 250777.75
$ws.Range("J126").Value = 2334.1667
$ws.Range("K126").Value = 752333.25
$ws.Range("L126").Value = 7002.500100000001
$ws.Range("M126").Value = -749863.25
$ws.Range("N126").Value = -11942.5001

$ws.Range("H134").Value = 134842.9
$ws.Range("J134").Value = 134842.9
$ws.Range("L134").Value = 404528.7
$ws.Range("N134").Value = -409598.7

$ws.Range("H136").Value = 19668602
$ws.Range("I136").Value = 27861842
$ws.Range("K136").Value = 83585526
$ws.Range("M136").Value = -83582976
